# Generate Report for Handoff
# Refresh the handoff report with a new source-file GUID, a new target-file
# hash and updated handoff timestamps.

$wb = $excel.ActiveWorkbook

$oldMdName  = "4a2ee02d-e25f-4e08-81a7-1621eb273dc1.md"
$newMdName  = "507f06c1-9115-4909-b3d1-b22155854e3e.md"

$oldZhXlf = "4a2ee02d-e25f-4e08-81a7-1621eb273dc1.e90acf5eddcc6550874b40128fb212b4c8232c8e.zh-cn.xlf"
$newZhXlf = "507f06c1-9115-4909-b3d1-b22155854e3e.aa8cd4b96cbdcb61133ce62069fa457cb6336f9e.zh-cn.xlf"

$oldDeXlf = "4a2ee02d-e25f-4e08-81a7-1621eb273dc1.e90acf5eddcc6550874b40128fb212b4c8232c8e.de-de.xlf"
$newDeXlf = "507f06c1-9115-4909-b3d1-b22155854e3e.aa8cd4b96cbdcb61133ce62069fa457cb6336f9e.de-de.xlf"

$newOverviewDate = "2016-46-17 06:46:44"
$newZhDate        = "2016-03-17 06:46:37"
$newDeDate        = "2016-03-17 06:46:44"

# ---------------------------------------------------------------------------
# Overview sheet: source file name (A2) + latest handoff date (D2)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newMdName
$wsOverview.Range("D2").Value = $newOverviewDate

foreach ($h in $wsOverview.Hyperlinks) {
    if ($h.Range().Address() -eq '$A$2') {
        $h.TextToDisplay = $newMdName
    }
}

# ---------------------------------------------------------------------------
# zh-cn sheet: source file name (A2), target file name (D2) + its date (E2)
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = $newMdName
$wsZhCn.Range("D2").Value = $newZhXlf
$wsZhCn.Range("E2").Value = $newZhDate

foreach ($h in $wsZhCn.Hyperlinks) {
    if ($h.Range().Address() -eq '$A$2') {
        $h.TextToDisplay = $newMdName
    } elseif ($h.Range().Address() -eq '$D$2') {
        $h.TextToDisplay = $newZhXlf
    }
}

# ---------------------------------------------------------------------------
# de-de sheet: source file name (A2), target file name (D2) + its date (E2)
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value = $newMdName
$wsDeDe.Range("D2").Value = $newDeXlf
$wsDeDe.Range("E2").Value = $newDeDate

foreach ($h in $wsDeDe.Hyperlinks) {
    if ($h.Range().Address() -eq '$A$2') {
        $h.TextToDisplay = $newMdName
    } elseif ($h.Range().Address() -eq '$D$2') {
        $h.TextToDisplay = $newDeXlf
    }
}
